# Re-theme the deck from the "Integral" (Red Violet) palette to the
# default Office Theme palette, and switch the three data tables from
# the plain "Table_0" style to the Accent-1 themed table style.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1. Theme colours: Integral / Red Violet  ->  Office Theme / Office
# ---------------------------------------------------------------------
function Set-ThemeRGB {
    param($Scheme, [int]$Index, [string]$Hex)
    $r = [Convert]::ToInt32($Hex.Substring(0,2), 16)
    $g = [Convert]::ToInt32($Hex.Substring(2,2), 16)
    $b = [Convert]::ToInt32($Hex.Substring(4,2), 16)
    $Scheme.Colors($Index).RGB = $r + ($g * 256) + ($b * 65536)
}

$tcs = $p.Slides.Item(1).ThemeColorScheme

Set-ThemeRGB $tcs 1  "000000"   # dk1
Set-ThemeRGB $tcs 2  "FFFFFF"   # lt1
Set-ThemeRGB $tcs 3  "44546A"   # dk2
Set-ThemeRGB $tcs 4  "E7E6E6"   # lt2
Set-ThemeRGB $tcs 5  "5B9BD5"   # accent1
Set-ThemeRGB $tcs 6  "ED7D31"   # accent2
Set-ThemeRGB $tcs 7  "A5A5A5"   # accent3
Set-ThemeRGB $tcs 8  "FFC000"   # accent4
Set-ThemeRGB $tcs 9  "4472C4"   # accent5
Set-ThemeRGB $tcs 10 "70AD47"   # accent6
Set-ThemeRGB $tcs 11 "0563C1"   # hlink
Set-ThemeRGB $tcs 12 "954F72"   # folHlink

# ---------------------------------------------------------------------
# 2. Table styles: plain "Table_0" -> themed "Medium Style 2 - Accent 1"
#    on every table in the deck (slides 14, 15 and 16 each hold one).
# ---------------------------------------------------------------------
$newTableStyle = "{BC5F0357-54D4-4088-9B34-3E82C7BF3446}"

for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)
    for ($j = 1; $j -le $slide.Shapes.Count; $j++) {
        $shape = $slide.Shapes.Item($j)
        if ($shape.HasTable) {
            $shape.Table.ApplyStyle($newTableStyle)
        }
    }
}
